# Applies the cryptos-list refresh described in the commit diff.
# Each entry is (cell reference, new value, forceText).
# forceText=$true means the literal string "looks like" a number to Excel
# (e.g. '0.9978' or '17.47'), so we briefly mark the cell as Text before
# writing it and restore the default "Normal" style afterwards so the
# saved cell keeps its original (unstyled) appearance but stores the exact
# literal digits as text, matching the source data feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "30.611.07"; ForceText = $false }
    @{ Cell = "E2"; Value = "  +0.95%  "; ForceText = $false }
    @{ Cell = "D3"; Value = "1.864.39"; ForceText = $false }
    @{ Cell = "E3"; Value = "  +0.12%  "; ForceText = $false }
    @{ Cell = "D4"; Value = "0.9978"; ForceText = $true }
    @{ Cell = "E4"; Value = "  -0.16%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "234.35"; ForceText = $true }
    @{ Cell = "E5"; Value = "  +0.15%  "; ForceText = $false }
    @{ Cell = "D6"; Value = "0.9977"; ForceText = $true }
    @{ Cell = "E6"; Value = "  -0.18%  "; ForceText = $false }
    @{ Cell = "D7"; Value = "0.4694"; ForceText = $true }
    @{ Cell = "E7"; Value = "  -1.53%  "; ForceText = $false }
    @{ Cell = "D8"; Value = "0.2765"; ForceText = $true }
    @{ Cell = "E8"; Value = "  +0.57%  "; ForceText = $false }
    @{ Cell = "D9"; Value = "0.06360"; ForceText = $true }
    @{ Cell = "E9"; Value = "  -1.37%  "; ForceText = $false }
    @{ Cell = "D10"; Value = "17.47"; ForceText = $true }
    @{ Cell = "E10"; Value = "  +8.26%  "; ForceText = $false }
    @{ Cell = "D11"; Value = "1.845.29"; ForceText = $false }
    @{ Cell = "E11"; Value = "  -0.22%  "; ForceText = $false }
    @{ Cell = "D12"; Value = "0.07460"; ForceText = $true }
    @{ Cell = "E12"; Value = "  +0.45%  "; ForceText = $false }
    @{ Cell = "D13"; Value = "4.965"; ForceText = $true }
    @{ Cell = "E13"; Value = "  -0.73%  "; ForceText = $false }
    @{ Cell = "D14"; Value = "84.97"; ForceText = $true }
    @{ Cell = "E14"; Value = "  -1.26%  "; ForceText = $false }
    @{ Cell = "D15"; Value = "0.6321"; ForceText = $true }
    @{ Cell = "E15"; Value = "  -0.11%  "; ForceText = $false }
    @{ Cell = "D16"; Value = "30.541.42"; ForceText = $false }
    @{ Cell = "E16"; Value = "  +0.77%  "; ForceText = $false }
    @{ Cell = "B17"; Value = "BitcoinCash"; ForceText = $false }
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; ForceText = $false }
    @{ Cell = "D17"; Value = "241.73"; ForceText = $true }
    @{ Cell = "E17"; Value = "  +3.74%  "; ForceText = $false }
    @{ Cell = "B18"; Value = "Dai"; ForceText = $false }
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; ForceText = $false }
    @{ Cell = "D18"; Value = "0.9984"; ForceText = $true }
    @{ Cell = "E18"; Value = "  -0.13%  "; ForceText = $false }
    @{ Cell = "D19"; Value = "12.72"; ForceText = $true }
    @{ Cell = "E19"; Value = "  -0.91%  "; ForceText = $false }
    @{ Cell = "D20"; Value = "0.000007362"; ForceText = $true }
    @{ Cell = "E20"; Value = "  -0.40%  "; ForceText = $false }
    @{ Cell = "D21"; Value = "0.9973"; ForceText = $true }
    @{ Cell = "E21"; Value = "  -0.16%  "; ForceText = $false }
    @{ Cell = "D22"; Value = "4.982"; ForceText = $true }
    @{ Cell = "E22"; Value = "  -2.56%  "; ForceText = $false }
    @{ Cell = "D23"; Value = "5.965"; ForceText = $true }
    @{ Cell = "E23"; Value = "  -1.05%  "; ForceText = $false }
    @{ Cell = "D24"; Value = "9.271"; ForceText = $true }
    @{ Cell = "E24"; Value = "  -0.40%  "; ForceText = $false }
    @{ Cell = "D25"; Value = "166.86"; ForceText = $true }
    @{ Cell = "E25"; Value = "  -0.57%  "; ForceText = $false }
    @{ Cell = "D26"; Value = "18.24"; ForceText = $true }
    @{ Cell = "E26"; Value = "  +1.79%  "; ForceText = $false }
    @{ Cell = "D27"; Value = "1.888"; ForceText = $true }
    @{ Cell = "E27"; Value = "  +1.43%  "; ForceText = $false }
    @{ Cell = "D28"; Value = "0.1026"; ForceText = $true }
    @{ Cell = "E28"; Value = "  +1.35%  "; ForceText = $false }
    @{ Cell = "D29"; Value = "1.378"; ForceText = $true }
    @{ Cell = "E29"; Value = "  -0.42%  "; ForceText = $false }
    @{ Cell = "D30"; Value = "4.118"; ForceText = $true }
    @{ Cell = "E30"; Value = "  -2.84%  "; ForceText = $false }
    @{ Cell = "D31"; Value = "3.863"; ForceText = $true }
    @{ Cell = "E31"; Value = "  -1.43%  "; ForceText = $false }
    @{ Cell = "D32"; Value = "0.04923"; ForceText = $true }
    @{ Cell = "E32"; Value = "  +0.35%  "; ForceText = $false }
    @{ Cell = "D33"; Value = "1.149"; ForceText = $true }
    @{ Cell = "E33"; Value = "  -0.30%  "; ForceText = $false }
    @{ Cell = "D34"; Value = "0.7085"; ForceText = $true }
    @{ Cell = "E34"; Value = "  -2.48%  "; ForceText = $false }
    @{ Cell = "D35"; Value = "2.697"; ForceText = $true }
    @{ Cell = "E35"; Value = "  +0.27%  "; ForceText = $false }
    @{ Cell = "D36"; Value = "0.01912"; ForceText = $true }
    @{ Cell = "E36"; Value = "  -1.66%  "; ForceText = $false }
    @{ Cell = "D37"; Value = "2.692"; ForceText = $true }
    @{ Cell = "E37"; Value = "  +2.32%  "; ForceText = $false }
    @{ Cell = "D38"; Value = "0.8834"; ForceText = $true }
    @{ Cell = "E38"; Value = "  -2.77%  "; ForceText = $false }
    @{ Cell = "D39"; Value = "1.974"; ForceText = $true }
    @{ Cell = "E39"; Value = "  -0.90%  "; ForceText = $false }
    @{ Cell = "D40"; Value = "105.84"; ForceText = $true }
    @{ Cell = "E40"; Value = "  +0.05%  "; ForceText = $false }
    @{ Cell = "D41"; Value = "0.9977"; ForceText = $true }
    @{ Cell = "E41"; Value = "  -0.21%  "; ForceText = $false }
    @{ Cell = "D42"; Value = "0.4097"; ForceText = $true }
    @{ Cell = "E42"; Value = "  -0.70%  "; ForceText = $false }
    @{ Cell = "D43"; Value = "5.545"; ForceText = $true }
    @{ Cell = "E43"; Value = "  -0.47%  "; ForceText = $false }
    @{ Cell = "D44"; Value = "7.214"; ForceText = $true }
    @{ Cell = "E44"; Value = "  +1.89%  "; ForceText = $false }
    @{ Cell = "D45"; Value = "0.1239"; ForceText = $true }
    @{ Cell = "E45"; Value = "  +2.49%  "; ForceText = $false }
    @{ Cell = "D46"; Value = "61.94"; ForceText = $true }
    @{ Cell = "E46"; Value = "  +0.83%  "; ForceText = $false }
    @{ Cell = "B47"; Value = "Elrond"; ForceText = $false }
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"; ForceText = $false }
    @{ Cell = "D47"; Value = "33.66"; ForceText = $true }
    @{ Cell = "E47"; Value = "  +1.75%  "; ForceText = $false }
    @{ Cell = "B48"; Value = "EnergySwap"; ForceText = $false }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; ForceText = $false }
    @{ Cell = "D48"; Value = "8.606"; ForceText = $true }
    @{ Cell = "E48"; Value = "  -1.81%  "; ForceText = $false }
    @{ Cell = "D49"; Value = "0.05553"; ForceText = $true }
    @{ Cell = "E49"; Value = "  -1.05%  "; ForceText = $false }
    @{ Cell = "D50"; Value = "1.379"; ForceText = $true }
    @{ Cell = "E50"; Value = "  -1.89%  "; ForceText = $false }
    @{ Cell = "D51"; Value = "0.3705"; ForceText = $true }
    @{ Cell = "E51"; Value = "  -0.42%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $r.NumberFormat = "@"
        $r.Value2 = $u.Value
        $r.Style = "Normal"
    } else {
        $r.Value2 = $u.Value
    }
}
